$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values, repulled/recalculated from source data
$updates = @{
    2  = 2
    8  = -3
    9  = -1
    10 = 1
    11 = -1
    17 = 3
    20 = -1
    21 = 0
    23 = -2
    24 = 0
    25 = -3
    27 = -1
    28 = 0
    31 = -3
    37 = -3
    38 = 1
    42 = -3
    46 = -2
    51 = 0
    53 = -7
    66 = 4
    67 = -5
    69 = -1
    71 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
